$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D retains its original text formatting so numeric-looking
# price strings (e.g. '1.031', '27.894.99') are not reinterpreted as numbers.
$ws.Range('D2:D51').NumberFormat = '@'

$ws.Range('D2').Value = '27.894.99'
$ws.Range('E2').Value = '  +1.92%  '
$ws.Range('D3').Value = '1.872.18'
$ws.Range('E3').Value = '  +1.92%  '
$ws.Range('D4').Value = '1.031'
$ws.Range('E4').Value = '  +0.42%  '
$ws.Range('D5').Value = '324.95'
$ws.Range('E5').Value = '  +2.24%  '
$ws.Range('D6').Value = '1.030'
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').Value = '0.4452'
$ws.Range('E7').Value = '  +2.20%  '
$ws.Range('D8').Value = '0.3830'
$ws.Range('E8').Value = '  +3.04%  '
$ws.Range('D9').Value = '0.07478'
$ws.Range('E9').Value = '  +1.89%  '
$ws.Range('D10').Value = '0.8916'
$ws.Range('E10').Value = '  +2.25%  '
$ws.Range('D11').Value = '21.81'
$ws.Range('E11').Value = '  +2.17%  '
$ws.Range('D12').Value = '1.878.43'
$ws.Range('E12').Value = '  -3.29%  '
$ws.Range('D13').Value = '5.596'
$ws.Range('E13').Value = '  +2.36%  '
$ws.Range('D14').Value = '6.794'
$ws.Range('E14').Value = '  +1.65%  '
$ws.Range('D15').Value = '0.07211'
$ws.Range('E15').Value = '  +1.25%  '
$ws.Range('D16').Value = '85.42'
$ws.Range('E16').Value = '  +4.02%  '
$ws.Range('D17').Value = '1.035'
$ws.Range('E17').Value = '  +0.66%  '
$ws.Range('D18').Value = '0.000009149'
$ws.Range('E18').Value = '  +1.92%  '
$ws.Range('D19').Value = '1.031'
$ws.Range('E19').Value = '  +0.85%  '
$ws.Range('D20').Value = '15.63'
$ws.Range('E20').Value = '  +1.74%  '
$ws.Range('D21').Value = '27.909.85'
$ws.Range('E21').Value = '  +1.87%  '
$ws.Range('D22').Value = '5.337'
$ws.Range('E22').Value = '  +1.73%  '
$ws.Range('D23').Value = '11.35'
$ws.Range('D24').Value = '2.105.05'
$ws.Range('E24').Value = '  -1.51%  '
$ws.Range('D25').Value = '2.024'
$ws.Range('E25').Value = '  +6.75%  '
$ws.Range('D26').Value = '158.97'
$ws.Range('E26').Value = '  +1.49%  '
$ws.Range('E27').Value = '  +2.19%  '
$ws.Range('D28').Value = '5.430'
$ws.Range('E28').Value = '  +3.68%  '
$ws.Range('D29').Value = '1.997'
$ws.Range('E29').Value = '  +4.05%  '
$ws.Range('D30').Value = '118.57'
$ws.Range('E30').Value = '  +2.67%  '
$ws.Range('D31').Value = '0.09086'
$ws.Range('E31').Value = '  +0.62%  '
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').Value = '0.7856'
$ws.Range('E33').Value = '  +3.48%  '
$ws.Range('D34').Value = '4.621'
$ws.Range('E34').Value = '  +3.64%  '
$ws.Range('E35').Value = '  +5.66%  '
$ws.Range('D36').Value = '1.033'
$ws.Range('E36').Value = '  +0.73%  '
$ws.Range('D37').Value = '1.150'
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.01995'
$ws.Range('E38').Value = '  +2.06%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.05379'
$ws.Range('E39').Value = '  +2.58%  '
$ws.Range('D40').Value = '2.881'
$ws.Range('E40').Value = '  +3.40%  '
$ws.Range('D41').Value = '0.5240'
$ws.Range('E41').Value = '  +1.63%  '
$ws.Range('D42').Value = '0.1701'
$ws.Range('E42').Value = '  +2.51%  '
$ws.Range('D43').Value = '6.949'
$ws.Range('E43').Value = '  +6.27%  '
$ws.Range('E44').Value = '  +4.74%  '
$ws.Range('D45').Value = '112.66'
$ws.Range('E45').Value = '  +4.11%  '
$ws.Range('D46').Value = '10.80'
$ws.Range('E46').Value = '  +2.65%  '
$ws.Range('D47').Value = '0.06618'
$ws.Range('E47').Value = '  +5.14%  '
$ws.Range('E48').Value = '  +3.40%  '
$ws.Range('D49').Value = '1.034'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('D50').Value = '0.4767'
$ws.Range('E50').Value = '  +3.27%  '
$ws.Range('D51').Value = '1.938'
$ws.Range('E51').Value = '  +3.16%  '
